# "Added two new Mac-Addresses" - append two more rows of reg_center/user/
# machine mapping data, following the exact pattern of the existing rows
# (regcntr_id=10001, sequential usr_id/machine_id, lang_code="eng",
# is_active=TRUE, cr_by="superadmin", cr_dtimes="now()").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 31; RegCntrId = 10001; UsrId = 110030; MachineId = 10030 },
    @{ Row = 32; RegCntrId = 10001; UsrId = 110031; MachineId = 10031 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.RegCntrId
    $ws.Cells.Item($row, 2).Value = $r.UsrId
    $ws.Cells.Item($row, 3).Value = $r.MachineId
    $ws.Cells.Item($row, 4).Value = "eng"
    $ws.Cells.Item($row, 5).Value = $true
    $ws.Cells.Item($row, 6).Value = "superadmin"
    $ws.Cells.Item($row, 7).Value = "now()"
}

# Match the saved selection/scroll state left behind after entering the
# new rows (cursor resting on C29, sheet scrolled down near the bottom).
$ws.Range("C29").Select()
